$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# n row
$ws.Range("D4").Value = 92
$ws.Range("E4").Value = 109

# Age (years), mean (SD)
$ws.Range("E5").Value = "7.8 (6.2)"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0.014"

# Age group (years), n (%)
$ws.Range("E6").Value = "45 (41.7)"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "0.084"
$ws.Range("E7").Value = "63 (58.3)"

# Sex, n (%)
$ws.Range("E8").Value = "50 (45.9)"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "0.507"
$ws.Range("E9").Value = "59 (54.1)"

# Race or ethnic group, n (%)
$ws.Range("E10").Value = "79 (73.8)"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "0.734"
$ws.Range("E11").Value = "17 (15.9)"
$ws.Range("E12").Value = "1 (0.9)"
$ws.Range("E13").Value = "1 (0.9)"
$ws.Range("E14").Value = "9 (8.4)"

# Hispanic or Latino ethnic group, n (%)
$ws.Range("E15").Value = "16 (14.8)"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0.407"
$ws.Range("E16").Value = "92 (85.2)"

# MRD 1 Status, n (%)
$ws.Range("E17").Value = "47 (46.1)"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "0.103"
$ws.Range("E18").Value = "55 (53.9)"

# Leucocyte counts (10^9/L), n (%)
$ws.Range("E19").Value = "44 (40.7)"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "0.388"
$ws.Range("E20").Value = "64 (59.3)"

# BM Leukemic blasts (%), mean (SD)
$ws.Range("E21").Value = "63.7 (26.8)"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "0.035"

# Risk Group, n (%)
$ws.Range("D22").Value = "15 (16.3)"
$ws.Range("E22").Value = "36 (33.0)"
$ws.Range("D23").Value = "27 (29.3)"
$ws.Range("E23").Value = "60 (55.0)"
$ws.Range("D24").Value = "50 (54.3)"
$ws.Range("E24").Value = "13 (11.9)"

# Clinical Trial, n (%)
$ws.Range("D25").Value = "74 (80.4)"
$ws.Range("E25").Value = "85 (78.0)"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "0.801"
$ws.Range("D26").Value = "18 (19.6)"
$ws.Range("E26").Value = "24 (22.0)"

# FLT3 ITD, n (%)
$ws.Range("D27").Value = "12 (13.2)"
$ws.Range("E27").Value = "19 (17.6)"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0.511"
$ws.Range("D28").Value = "79 (86.8)"
$ws.Range("E28").Value = "89 (82.4)"
